$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.167.03'
$ws.Range("E2").Value = '  -10.68%  '
$ws.Range("D3").Value = '2.302.61'
$ws.Range("E3").Value = '  -20.69%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '447.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -15.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.71'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -10.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.472'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -14.95%  '
$ws.Range("D9").Value = '2.304.71'
$ws.Range("E9").Value = '  -20.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -10.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0917'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -15.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.309'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -14.88%  '
$ws.Range("E13").Value = '  -3.47%  '
$ws.Range("D14").Value = '2.696.16'
$ws.Range("E14").Value = '  -20.76%  '
$ws.Range("D15").Value = '54.162.57'
$ws.Range("E15").Value = '  -10.67%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.80'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -17.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000121'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -14.92%  '
$ws.Range("D18").Value = '2.308.78'
$ws.Range("E18").Value = '  -20.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -19.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '301.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -16.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -19.50%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -20.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '55.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -14.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.994'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.49%  '
$ws.Range("E27").Value = '  -12.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.370'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -18.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.82'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -13.49%  '
$ws.Range("D31").Value = '0.0₃0707'
$ws.Range("E31").Value = '  -17.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '145.14'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '16.94'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -14.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.35'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -19.82%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -15.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -17.67%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.843'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -16.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.01'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -16.16%  '
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '33.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.08%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("E42").Value = '  -16.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -15.46%  '
$ws.Range("D44").Value = '1.931.43'
$ws.Range("E44").Value = '  -15.78%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0497'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -14.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.514'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -20.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0206'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -13.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0815'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -11.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -21.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -19.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.82%  '
